$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '28.373.08'
$r.ClearFormats()
$ws.Range("E2").Value = '  +5.29%  '

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '1.799.88'
$r.ClearFormats()
$ws.Range("E3").Value = '  +3.39%  '

$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '1.003'
$r.ClearFormats()
$ws.Range("E4").Value = '  +0.18%  '

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '315.76'
$r.ClearFormats()
$ws.Range("E5").Value = '  +1.19%  '

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '1.002'
$r.ClearFormats()
$ws.Range("E6").Value = '  +0.21%  '

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.5463'
$r.ClearFormats()
$ws.Range("E7").Value = '  +8.95%  '

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.3822'
$r.ClearFormats()
$ws.Range("E8").Value = '  +6.92%  '

$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '42.91'
$r.ClearFormats()
$ws.Range("E9").Value = '  +0.82%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.07547'
$r.ClearFormats()
$ws.Range("E10").Value = '  +3.91%  '

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '1.120'
$r.ClearFormats()
$ws.Range("E11").Value = '  +5.75%  '

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '1.002'
$r.ClearFormats()
$ws.Range("E12").Value = '  +0.17%  '

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '21.09'
$r.ClearFormats()
$ws.Range("E13").Value = '  +3.76%  '

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '6.175'
$r.ClearFormats()
$ws.Range("E14").Value = '  +3.08%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '1.799.59'
$r.ClearFormats()
$ws.Range("E15").Value = '  +3.22%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '7.273'
$r.ClearFormats()
$ws.Range("E16").Value = '  +6.01%  '

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '90.92'
$r.ClearFormats()
$ws.Range("E17").Value = '  +4.64%  '

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '0.00001065'
$r.ClearFormats()
$ws.Range("E18").Value = '  +2.98%  '

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '0.06455'
$r.ClearFormats()
$ws.Range("E19").Value = '  +0.91%  '

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '1.002'
$r.ClearFormats()
$ws.Range("E20").Value = '  +0.20%  '

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '17.20'
$r.ClearFormats()
$ws.Range("E21").Value = '  +3.75%  '

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '5.960'
$r.ClearFormats()
$ws.Range("E22").Value = '  +3.97%  '

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '28.383.33'
$r.ClearFormats()
$ws.Range("E23").Value = '  +4.97%  '

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '11.19'
$r.ClearFormats()
$ws.Range("E24").Value = '  -0.23%  '

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '2.115'
$r.ClearFormats()
$ws.Range("E25").Value = '  +3.24%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '158.36'
$r.ClearFormats()
$ws.Range("E26").Value = '  +2.64%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '20.59'
$r.ClearFormats()
$ws.Range("E27").Value = '  +2.88%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '2.410'
$r.ClearFormats()
$ws.Range("E28").Value = '  +12.43%  '

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '2.008.05'
$r.ClearFormats()
$ws.Range("E29").Value = '  +3.25%  '

$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '123.62'
$r.ClearFormats()
$ws.Range("E30").Value = '  +1.81%  '

$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '1.154'
$r.ClearFormats()
$ws.Range("E31").Value = '  +9.94%  '

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '0.1024'
$r.ClearFormats()
$ws.Range("E32").Value = '  +7.63%  '

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '5.690'
$r.ClearFormats()
$ws.Range("E33").Value = '  +5.48%  '

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '3.674'
$r.ClearFormats()
$ws.Range("E34").Value = '  +2.85%  '

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '0.2257'
$r.ClearFormats()
$ws.Range("E35").Value = '  +12.91%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '0.06289'
$r.ClearFormats()
$ws.Range("E36").Value = '  +6.67%  '

$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '8.853'
$r.ClearFormats()
$ws.Range("E37").Value = '  +17.69%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '0.02308'
$r.ClearFormats()
$ws.Range("E38").Value = '  +4.36%  '

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '11.52'
$r.ClearFormats()
$ws.Range("E39").Value = '  +3.77%  '

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '4.992'
$r.ClearFormats()
$ws.Range("E40").Value = '  +4.95%  '

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.6336'
$r.ClearFormats()
$ws.Range("E41").Value = '  +5.04%  '

$ws.Range("B42").Value = 'Frax'
$ws.Range("C42").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '1.002'
$r.ClearFormats()
$ws.Range("E42").Value = '  +0.20%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '1.153'
$r.ClearFormats()
$ws.Range("E43").Value = '  +3.08%  '

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '1.382'
$r.ClearFormats()
$ws.Range("E44").Value = '  -3.63%  '

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '13.41'
$r.ClearFormats()
$ws.Range("E45").Value = '  +4.55%  '

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '0.5930'
$r.ClearFormats()
$ws.Range("E46").Value = '  +4.95%  '

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '3.668'
$r.ClearFormats()
$ws.Range("E47").Value = '  +1.84%  '

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '123.28'
$r.ClearFormats()
$ws.Range("E48").Value = '  +2.74%  '

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '1.953'
$r.ClearFormats()
$ws.Range("E49").Value = '  +4.90%  '

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '1.142'
$r.ClearFormats()
$ws.Range("E50").Value = '  +3.04%  '

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '0.06897'
$r.ClearFormats()
$ws.Range("E51").Value = '  +3.38%  '
